# Update odds values in the "Jogos da Semana" worksheet to reflect the
# latest FlashScore data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Aston Villa vs Liverpool)
$ws.Range("O2").Value = 1.17
$ws.Range("P2").Value = 5
$ws.Range("Q2").Value = 1.57
$ws.Range("R2").Value = 2.38
$ws.Range("S2").Value = 2
$ws.Range("T2").Value = 1.9

# Row 3 (Once Caldas vs Pereira)
$ws.Range("M3").Value = 1.08
$ws.Range("O3").Value = 1.44
$ws.Range("P3").Value = 2.63
$ws.Range("V3").Value = 1.18

# Row 4 (Luton vs Plymouth)
$ws.Range("M4").Value = 1.08
$ws.Range("O4").Value = 1.4
$ws.Range("U4").Value = 4.33

# Row 5 (Hyderabad vs Mumbai City)
$ws.Range("H5").Value = 4
$ws.Range("I5").Value = 1.6
$ws.Range("K5").Value = 2.4
$ws.Range("L5").Value = 2.2
$ws.Range("N5").Value = 15
$ws.Range("O5").Value = 1.18
$ws.Range("P5").Value = 4.5
$ws.Range("Q5").Value = 1.62
$ws.Range("R5").Value = 2.25
$ws.Range("S5").Value = 2
$ws.Range("T5").Value = 1.8
$ws.Range("U5").Value = 2.5
$ws.Range("V5").Value = 1.5
$ws.Range("W5").Value = 1.3
$ws.Range("X5").Value = 3.4
$ws.Range("Y5").Value = 1.67
$ws.Range("Z5").Value = 2.1
$ws.Range("AA5").Value = 17
$ws.Range("AB5").Value = 29
$ws.Range("AE5").Value = 34
$ws.Range("AF5").Value = 34
$ws.Range("AG5").Value = 15
$ws.Range("AH5").Value = 8
$ws.Range("AJ5").Value = 41
$ws.Range("AK5").Value = 151
$ws.Range("AL5").Value = 9
$ws.Range("AM5").Value = 9
$ws.Range("AP5").Value = 12
$ws.Range("AQ5").Value = 21

# Row 7 (Club America vs Club Leon)
$ws.Range("G7").Value = 1.57
$ws.Range("H7").Value = 4.33
$ws.Range("I7").Value = 5.25
$ws.Range("J7").Value = 2.1
$ws.Range("N7").Value = 13
$ws.Range("Q7").Value = 1.73
$ws.Range("R7").Value = 2.08
$ws.Range("AB7").Value = 7.5
$ws.Range("AG7").Value = 13
$ws.Range("AH7").Value = 8
$ws.Range("AI7").Value = 17
$ws.Range("AK7").Value = 251
